$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 35719388
$ws.Range("I70").Value = 2724.75
$ws.Range("K70").Value = 8174.25
$ws.Range("M70").Value = -7904.25
$ws.Range("H73").Value = 35719388
$ws.Range("I73").Value = 2724.75
$ws.Range("K73").Value = 8174.25
$ws.Range("M73").Value = -7238.25
$ws.Range("H137").Value = 11446202
$ws.Range("I137").Value = 401722.56
$ws.Range("J137").Value = 66668596
$ws.Range("K137").Value = 1205167.68
$ws.Range("L137").Value = 200005788
$ws.Range("M137").Value = -1202617.68
$ws.Range("N137").Value = -200010888
$ws.Range("H138").Value = 4328.811
$ws.Range("J138").Value = 6845.617
$ws.Range("L138").Value = 20536.851
$ws.Range("N138").Value = -30816.851

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15013.523
$ws.Range("I32").Value = 15998.412
$ws.Range("J32").Value = 11425.714
$ws.Range("K32").Value = 15998.412
$ws.Range("L32").Value = 11425.714
$ws.Range("M32").Value = -15711.412
$ws.Range("N32").Value = -11999.714
$ws.Range("H74").Value = 1112.3489
$ws.Range("I74").Value = 700.24243
$ws.Range("J74").Value = 2472.3
$ws.Range("K74").Value = 700.24243
$ws.Range("L74").Value = 2472.3
$ws.Range("M74").Value = 173.75757
$ws.Range("N74").Value = -4220.3
$ws.Range("H77").Value = 1112.3489
$ws.Range("I77").Value = 700.24243
$ws.Range("J77").Value = 2472.3
$ws.Range("K77").Value = 3501.21215
$ws.Range("L77").Value = 12361.5
$ws.Range("M77").Value = 866.7878499999997
$ws.Range("N77").Value = -21097.5
$ws.Range("H124").Value = 61000
$ws.Range("J124").Value = 61000
$ws.Range("L124").Value = 61000
$ws.Range("N124").Value = -70820
$ws.Range("H132").Value = 25505.178
$ws.Range("I132").Value = 29770.8
$ws.Range("K132").Value = 89312.39999999999
$ws.Range("M132").Value = -86782.39999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26318700
$ws.Range("I31").Value = 31252386
$ws.Range("J31").Value = 5709.1665
$ws.Range("K31").Value = 31252386
$ws.Range("L31").Value = 5709.1665
$ws.Range("M31").Value = -31252091
$ws.Range("N31").Value = -6299.1665
$ws.Range("H34").Value = 26318700
$ws.Range("I34").Value = 31252386
$ws.Range("J34").Value = 5709.1665
$ws.Range("K34").Value = 31252386
$ws.Range("L34").Value = 5709.1665
$ws.Range("M34").Value = -31252184
$ws.Range("N34").Value = -6113.1665
$ws.Range("H122").Value = 2820.9697
$ws.Range("I122").Value = 1455.1305
$ws.Range("K122").Value = 4365.3915
$ws.Range("M122").Value = -1915.3915
$ws.Range("H132").Value = 28377356
$ws.Range("I132").Value = 41674910
$ws.Range("J132").Value = 9231.866
$ws.Range("K132").Value = 125024730
$ws.Range("L132").Value = 27695.598
$ws.Range("M132").Value = -125022200
$ws.Range("N132").Value = -32755.598
$ws.Range("H134").Value = 2041.4849
$ws.Range("I134").Value = 2009.3334
$ws.Range("J134").Value = 2186.1667
$ws.Range("K134").Value = 6028.0002
$ws.Range("L134").Value = 6558.500100000001
$ws.Range("M134").Value = -3493.0002
$ws.Range("N134").Value = -11628.5001
$ws.Range("H141").Value = 109200.91
$ws.Range("J141").Value = 117590.11
$ws.Range("L141").Value = 117590.11
$ws.Range("N141").Value = -127950.11

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 5904.7144
$ws.Range("I82").Value = 3333
$ws.Range("J82").Value = 6333.3335
$ws.Range("K82").Value = 9999
$ws.Range("L82").Value = 19000.0005
$ws.Range("M82").Value = -9593
$ws.Range("N82").Value = -19812.0005
$ws.Range("H85").Value = 5904.7144
$ws.Range("I85").Value = 3333
$ws.Range("J85").Value = 6333.3335
$ws.Range("K85").Value = 9999
$ws.Range("L85").Value = 19000.0005
$ws.Range("M85").Value = -8595
$ws.Range("N85").Value = -21808.0005
$ws.Range("H132").Value = 4902.706
$ws.Range("I132").Value = 1229.1818
$ws.Range("J132").Value = 11637.5
$ws.Range("K132").Value = 11062.6362
$ws.Range("L132").Value = 104737.5
$ws.Range("M132").Value = -8532.636200000001
$ws.Range("N132").Value = -109797.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1307424.2
$ws.Range("I80").Value = 2238129.8
$ws.Range("J80").Value = 4436.4
$ws.Range("K80").Value = 2238129.8
$ws.Range("L80").Value = 4436.4
$ws.Range("M80").Value = -2237131.8
$ws.Range("N80").Value = -6432.4
$ws.Range("H83").Value = 1307424.2
$ws.Range("I83").Value = 2238129.8
$ws.Range("J83").Value = 4436.4
$ws.Range("K83").Value = 11190649
$ws.Range("L83").Value = 22182
$ws.Range("M83").Value = -11185657
$ws.Range("N83").Value = -32166
$ws.Range("H113").Value = 3688.3
$ws.Range("I113").Value = 3416.8
$ws.Range("J113").Value = 3959.8
$ws.Range("K113").Value = 3416.8
$ws.Range("L113").Value = 3959.8
$ws.Range("M113").Value = -1246.8
$ws.Range("N113").Value = -8299.799999999999
$ws.Range("H122").Value = 398799.94
$ws.Range("I122").Value = 918404.5600000001
$ws.Range("J122").Value = 9096.4375
$ws.Range("K122").Value = 2755213.68
$ws.Range("L122").Value = 27289.3125
$ws.Range("M122").Value = -2752763.68
$ws.Range("N122").Value = -32189.3125
$ws.Range("H123").Value = 56005.5
$ws.Range("J123").Value = 56005.5
$ws.Range("L123").Value = 56005.5
$ws.Range("N123").Value = -60905.5
$ws.Range("H126").Value = 3606.7
$ws.Range("I126").Value = 2248.08
$ws.Range("J126").Value = 10399.8
$ws.Range("K126").Value = 6744.24
$ws.Range("L126").Value = 31199.4
$ws.Range("M126").Value = -4274.24
$ws.Range("N126").Value = -36139.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1213.4762
$ws.Range("I22").Value = 968.9
$ws.Range("K22").Value = 968.9
$ws.Range("M22").Value = -673.9
$ws.Range("H27").Value = 1213.4762
$ws.Range("I27").Value = 968.9
$ws.Range("K27").Value = 968.9
$ws.Range("M27").Value = -861.9
$ws.Range("H40").Value = 15629974
$ws.Range("I40").Value = 17862112
$ws.Range("K40").Value = 17862112
$ws.Range("M40").Value = -17861976
$ws.Range("H132").Value = 3281.6511
$ws.Range("J132").Value = 3550.5
$ws.Range("L132").Value = 10651.5
$ws.Range("N132").Value = -15711.5
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 660.41174
$ws.Range("I107").Value = 587.5
$ws.Range("J107").Value = 1000.6667
$ws.Range("K107").Value = 1762.5
$ws.Range("L107").Value = 3002.0001
$ws.Range("M107").Value = 157.5
$ws.Range("N107").Value = -6842.0001
